# Mise à jour 25/06
# - Add a new column F "statut_espece_lrr_2015" with red-list status codes
#   (CR / NT / LC / DD) for the native ("indigene") species rows.
# - Apply an AutoFilter on column D (statut_esp) to show only "indigene"
#   rows, which hides the non-native rows (2-22).
# - Correct E27 (temps_generation for Anguille européenne) from 12 to 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header ---------------------------------------------------
$ws.Range("F1").Value = "statut_espece_lrr_2015"

# --- Data correction -------------------------------------------------------
$ws.Range("E27").Value = 11

# --- New column F values for the relevant (indigene) rows -----------------
$ws.Range("F27").Value = "CR"
$ws.Range("F30").Value = "NT"
$ws.Range("F32").Value = "LC"
$ws.Range("F33").Value = "LC"
$ws.Range("F37").Value = "LC"
$ws.Range("F38").Value = "LC"
$ws.Range("F39").Value = "LC"
$ws.Range("F41").Value = "LC"
$ws.Range("F44").Value = "LC"
$ws.Range("F46").Value = "NT"
$ws.Range("F48").Value = "LC"
$ws.Range("F49").Value = "LC"
$ws.Range("F50").Value = "DD"

# --- Column width for the new column ---------------------------------------
$ws.Columns.Item(6).ColumnWidth = 23

# --- Filter the table on column D ("statut_esp") to just "indigene" --------
# colId 3 = 4th column of the A1:F50 range (D), xlFilterValues = 7
$ws.Range("A1:F50").AutoFilter(4, "indigene", 7) | Out-Null

# --- Keep the hidden _FilterDatabase defined name in sync with the new range
$filterName = $wb.Names.Item(1)
$filterName.RefersTo = "=Feuil1!`$A`$1:`$F`$50"

# --- Update the selection to match the author's last position --------------
$ws.Range("H25").Select() | Out-Null
